$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 4.692200000000001
$ws.Range("B10").Value = 8.590000000000005
$ws.Range("B12").Value = 6.6606
$ws.Range("D13").Value = -7.441100000000005
$ws.Range("B18").Value = 6.528000000000007
